# Update the division-practice table: each data row (1, 5, 9, 13, 17)
# holds 5 "a÷b=c, r" problems across the 5 columns. Replace each cell's
# text in place (by row/column position) so the run formatting
# (TimeNewRoman, sz 30) carries over untouched, and so the fact that a
# couple of the new values coincide with other cells' old values never
# causes an ambiguous/duplicate Find match.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "29÷4=7, 1"
$t.Cell(1, 2).Range.Text = "35÷4=8, 3"
$t.Cell(1, 3).Range.Text = "46÷5=9, 1"
$t.Cell(1, 4).Range.Text = "56÷7=8, 0"
$t.Cell(1, 5).Range.Text = "62÷4=15, 2"

$t.Cell(5, 1).Range.Text = "72÷3=24, 0"
$t.Cell(5, 2).Range.Text = "41÷2=20, 1"
$t.Cell(5, 3).Range.Text = "13÷4=3, 1"
$t.Cell(5, 4).Range.Text = "22÷6=3, 4"
$t.Cell(5, 5).Range.Text = "80÷9=8, 8"

$t.Cell(9, 1).Range.Text = "99÷2=49, 1"
$t.Cell(9, 2).Range.Text = "74÷5=14, 4"
$t.Cell(9, 3).Range.Text = "32÷3=10, 2"
$t.Cell(9, 4).Range.Text = "22÷3=7, 1"
$t.Cell(9, 5).Range.Text = "41÷9=4, 5"

$t.Cell(13, 1).Range.Text = "10÷3=3, 1"
$t.Cell(13, 2).Range.Text = "52÷5=10, 2"
$t.Cell(13, 3).Range.Text = "13÷4=3, 1"
$t.Cell(13, 4).Range.Text = "38÷7=5, 3"
$t.Cell(13, 5).Range.Text = "50÷7=7, 1"

$t.Cell(17, 1).Range.Text = "39÷9=4, 3"
$t.Cell(17, 2).Range.Text = "28÷8=3, 4"
$t.Cell(17, 3).Range.Text = "70÷2=35, 0"
$t.Cell(17, 4).Range.Text = "15÷3=5, 0"
$t.Cell(17, 5).Range.Text = "20÷7=2, 6"
